$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45882
$ws.Range("B2").Value = 109.73
$ws.Range("C2").Value = 110.7
$ws.Range("D2").Value = 109.98
$ws.Range("E2").Value = 105.55
$ws.Range("F2").Value = 104.32
$ws.Range("G2").Value = 104.32
$ws.Range("H2").Value = 105
$ws.Range("I2").Value = 121.17
$ws.Range("J2").Value = 118.24
$ws.Range("K2").Value = 94.09999999999999
$ws.Range("L2").Value = 90
$ws.Range("M2").Value = 71
$ws.Range("N2").Value = 63.48
$ws.Range("O2").Value = 60
$ws.Range("P2").Value = 63.22
$ws.Range("Q2").Value = 70.90000000000001
$ws.Range("R2").Value = 84.15000000000001
$ws.Range("S2").Value = 93.09999999999999
$ws.Range("T2").Value = 105
$ws.Range("U2").Value = 107.23
$ws.Range("V2").Value = 155.1
$ws.Range("W2").Value = 153
$ws.Range("X2").Value = 123.21
$ws.Range("Y2").Value = 106.27
$ws.Range("Z2").Value = 101.2
$ws.Range("AB2").Value = 134.39
$ws.Range("AD2").Value = 154.05
$ws.Range("AF2").Value = 114.74
$ws.Range("AG2").Value = "9h-17h"
